$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 65248
$ws.Cells.Item(2, 2).Value = "Giovanna Marques"
$ws.Cells.Item(2, 5).Value = 4
$ws.Cells.Item(2, 6).Value = 45100
$ws.Cells.Item(2, 7).Value = 8718.209999999999

# Row 3
$ws.Cells.Item(3, 1).Value = 19486
$ws.Cells.Item(3, 2).Value = "Benjamin Fonseca"
$ws.Cells.Item(3, 3).Value = "Engenharia"
$ws.Cells.Item(3, 4).Value = "Viagem de negocios"
$ws.Cells.Item(3, 6).Value = 45087
$ws.Cells.Item(3, 7).Value = 2737.28

# Row 4
$ws.Cells.Item(4, 1).Value = 71066
$ws.Cells.Item(4, 2).Value = "Clara Mendes"
$ws.Cells.Item(4, 3).Value = "Operacoes"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 45080
$ws.Cells.Item(4, 7).Value = 2812.6

# Row 5
$ws.Cells.Item(5, 1).Value = 42769
$ws.Cells.Item(5, 2).Value = "Maria Cecília Nunes"
$ws.Cells.Item(5, 3).Value = "TI"
$ws.Cells.Item(5, 4).Value = "Outros"
$ws.Cells.Item(5, 6).Value = 45097
$ws.Cells.Item(5, 7).Value = 6991.95

# Row 6
$ws.Cells.Item(6, 1).Value = 93984
$ws.Cells.Item(6, 2).Value = "Felipe Peixoto"
$ws.Cells.Item(6, 4).Value = "Problemas pessoais"
$ws.Cells.Item(6, 5).Value = 7
$ws.Cells.Item(6, 6).Value = 45105
$ws.Cells.Item(6, 7).Value = 4060.42

# Row 7
$ws.Cells.Item(7, 1).Value = 99449
$ws.Cells.Item(7, 2).Value = "Sr. Thales Cavalcanti"
$ws.Cells.Item(7, 3).Value = "Operacoes"
$ws.Cells.Item(7, 5).Value = 6
$ws.Cells.Item(7, 6).Value = 45102
$ws.Cells.Item(7, 7).Value = 3656

# Row 8
$ws.Cells.Item(8, 1).Value = 29718
$ws.Cells.Item(8, 2).Value = "Srta. Catarina Vieira"
$ws.Cells.Item(8, 3).Value = "Operacoes"
$ws.Cells.Item(8, 4).Value = "Viagem de negocios"
$ws.Cells.Item(8, 5).Value = 4
$ws.Cells.Item(8, 6).Value = 45079
$ws.Cells.Item(8, 7).Value = 5841.91

# Row 9
$ws.Cells.Item(9, 1).Value = 10403
$ws.Cells.Item(9, 2).Value = "Levi Cunha"
$ws.Cells.Item(9, 3).Value = "Vendas"
$ws.Cells.Item(9, 4).Value = "Doenca"
$ws.Cells.Item(9, 5).Value = 4
$ws.Cells.Item(9, 6).Value = 45091
$ws.Cells.Item(9, 7).Value = 5911.14

# Row 10
$ws.Cells.Item(10, 1).Value = 35646
$ws.Cells.Item(10, 2).Value = "Brenda Ramos"
$ws.Cells.Item(10, 3).Value = "Engenharia"
$ws.Cells.Item(10, 4).Value = "Viagem de negocios"
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = 45085
$ws.Cells.Item(10, 7).Value = 7396.91

# Row 11
$ws.Cells.Item(11, 1).Value = 33075
$ws.Cells.Item(11, 2).Value = "Diego Lopes"
$ws.Cells.Item(11, 5).Value = 5
$ws.Cells.Item(11, 6).Value = 45104
$ws.Cells.Item(11, 7).Value = 2889.22
